# Generate Report for Handback
#
# For the "ae299e66-1856-40c5-bd79-c1846a873e5a" entry (row 7 in both the
# zh-cn and de-de sheets) the handback was processed: the latest handback
# file/date are now known but the handback turned out to be based on a
# stale handoff version, so an error message is recorded in column P
# ("Error Detail").

$wb = $excel.ActiveWorkbook

$warningMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ae643981e95488744450c83c7262506db2794833/e2e/ae299e66-1856-40c5-bd79-c1846a873e5a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6ffba78403438fff811cc42f03fe982897aade2b/e2e/ae299e66-1856-40c5-bd79-c1846a873e5a.md."

# ---------------------------------------------------------------------
# zh-cn sheet, row 7
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# I7: Latest Target File -> becomes a hyperlink to the source markdown
# file (matches the pattern already used for rows 2,3,4,5 in column I).
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/6ffba78403438fff811cc42f03fe982897aade2b/e2e/ae299e66-1856-40c5-bd79-c1846a873e5a.md", "", "", "ae299e66-1856-40c5-bd79-c1846a873e5a.md") | Out-Null

# J7: Latest Handback File -> same xlf file name already shown in G7.
$wsZhCn.Range("J7").Value = $wsZhCn.Range("G7").Value2

# K7: Latest Handback DateTime
$wsZhCn.Range("K7").Value = "2016-08-28 12:56:01"

# P7: Error Detail
$wsZhCn.Range("P7").Value = $warningMessage

# ---------------------------------------------------------------------
# de-de sheet, row 7
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# I7: Latest Target File -> hyperlink, same as zh-cn.
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6ffba78403438fff811cc42f03fe982897aade2b/e2e/ae299e66-1856-40c5-bd79-c1846a873e5a.md", "", "", "ae299e66-1856-40c5-bd79-c1846a873e5a.md") | Out-Null

# J7: Latest Handback File -> same xlf file name already shown in G7.
$wsDeDe.Range("J7").Value = $wsDeDe.Range("G7").Value2

# K7: Latest Handback DateTime
$wsDeDe.Range("K7").Value = "2016-08-28 12:56:11"

# P7: Error Detail (identical wording to zh-cn row)
$wsDeDe.Range("P7").Value = $warningMessage
